$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "https://dashanddata.com"
$ws.Range("C2").Value = "max-age=63072000; includeSubdomains"
